$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Sprint 2 table data updates (rows 17-30)
# ---------------------------------------------------------------------------

# E17: estimated points for Sprint 2 changes from 34 to 22
$ws.Range("E17").Value = 22

# G17 / G18 formulas now reference $E$17 (Sprint-2 estimate) instead of $E$3
$ws.Range("G17").Formula = "=`$E`$17-F17"
$ws.Range("G18:G30").Formula = "=`$E`$17-F18"

# C18: second-series start date changes from 43772 to 43802
$ws.Range("C18").Value = 43802

# New "Points left" burn-down series for the Sprint-2 chart (D21:E23)
$ws.Range("D21").Value = 43789
$ws.Range("D22").Value = 43799
$ws.Range("D23").Value = 43801

$ws.Range("E21").Formula = "=G17"
$ws.Range("E22").Formula = "=G27"
$ws.Range("E23").Formula = "=G29"

# Copy number/date formatting onto the newly populated cells so the styles
# match the rest of the table (date style for column D, number style for E)
$ws.Range("C17").Copy()
$ws.Range("D21:D23").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("E21:E23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Daily points-completed entries that were previously blank
$ws.Range("B28").Value = 0
$ws.Range("B29").Value = 5

# ---------------------------------------------------------------------------
# 2. Selection / view bookkeeping
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 3
$ws.Range("U26").Select()

# ---------------------------------------------------------------------------
# 3. New "Sprint 2" burn-down chart (mirrors the existing Sprint-1 chart)
# ---------------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Add(700, 300, 509, 268)
$chartObj.Name = "Grafico 4"
$chart = $chartObj.Chart
$chart.ChartType = 74

$s1 = $chart.SeriesCollection().NewSeries()
$s1.Name = "=Foglio1!`$G`$2"
$s1.XValues = $ws.Range("D21:D23")
$s1.Values = $ws.Range("E21:E23")
$s1.MarkerStyle = 8
$s1.Smooth = $false

$s2 = $chart.SeriesCollection().NewSeries()
$s2.Name = "=Foglio1!`$H`$2"
$s2.XValues = $ws.Range("C17:C18")
$s2.Values = $ws.Range("H17:H18")
$s2.MarkerStyle = 8
$s2.Smooth = $false

$chart.HasTitle = $false
